# Applies the "Online Course Resources" section rewrite:
#  - "...course platform" -> "...course platform:" (colon added)
#  - trims / rewrites the trailing sentence about Blackboard
#  - adds two new BodyText paragraphs: one pointing at the "course drive"
#    and one restating the attendance note (replacing text removed from
#    the paragraph above)

$d = $word.ActiveDocument

# 1) Shorten / reword the long explanatory sentence, and use the freed-up
#    paragraph break to grow two new BodyText paragraphs in its place.
#    (A unique placeholder marks where the new hyperlink belongs so it can
#    be turned into a real w:hyperlink afterwards.)
$oldTail = ". In it, you will find submission portals for (some) of your assignments and a link to this course webpage, where you can find the course syllabus, problem sets, and links to readings. In addition, during the semester, solutions to the problem sets and lecture slides will be posted. Please note, however, that class attendance is the primary source of course-related announcements and material."
$newTail = ". In it, you will find submission portals for assignments and a link to this course webpage, where you can find the course syllabus, problem sets, and links to readings.^pYou can also link to our COURSEDRIVEPLACEHOLDER, which contains lecture slides, data sets, and some other useful things for the class."

[void]$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)

# 2) Add the colon after "...course platform" (kept as its own step so the
#    match stays unambiguous / scoped to that exact run).
$f1 = $d.Content.Find
$f1.Text = "Blackboard is our internet-based course platform"
$f1.Forward = $true
$f1.Wrap = 1
$f1.MatchCase = $true
[void]$f1.Execute()
$f1.Parent.Text = "Blackboard is our internet-based course platform:"

# 3) Split off the final reminder sentence into its own new paragraph
#    (kept as a separate Find/Replace call -- doing both paragraph breaks
#    in one call mangles the character formatting of the middle paragraph).
$oldNote = "which contains lecture slides, data sets, and some other useful things for the class."
$newNote = "which contains lecture slides, data sets, and some other useful things for the class.^pPlease note that class attendance is the primary source of course-related announcements and material."
[void]$d.Content.Find.Execute($oldNote, $true, $false, $false, $false, $false, $true, 1, $false, $newNote, 2)

# 4) The two new paragraphs inherited the "FirstParagraph" style from the
#    paragraph they split off of; re-style them to "BodyText" per the diff.
$f2 = $d.Content.Find
$f2.Text = "You can also link to our COURSEDRIVEPLACEHOLDER"
$f2.Forward = $true
$f2.Wrap = 1
$f2.MatchCase = $true
[void]$f2.Execute()
$linkPara = $f2.Parent.Paragraphs(1)
$linkPara.Style = "BodyText"
$notePara = $linkPara.Next()
$notePara.Style = "BodyText"

# 5) Turn the placeholder into a real hyperlink reading "course drive here".
$f3 = $d.Content.Find
$f3.Text = "COURSEDRIVEPLACEHOLDER"
$f3.Forward = $true
$f3.Wrap = 1
$f3.MatchCase = $true
[void]$f3.Execute()
$linkRange = $f3.Parent
[void]$d.Hyperlinks.Add($linkRange, "https://drive.google.com/", $null, $null, "course drive here")

Write-Output "done"
